$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Sex)
$ws.Range("F2").Value = 100
$ws.Range("H2").Value = 100
$ws.Range("I2").Value = 26.5323113860929
$ws.Range("J2").Value = 0
$ws.Range("K2").ClearContents()
$ws.Range("L2").Value = "Significant"

# Row 3 (Male)
$ws.Range("F3").Value = 94.12
$ws.Range("H3").Value = 21.43
$ws.Range("L3").ClearContents()

# Row 4 (Female)
$ws.Range("F4").Value = 5.88
$ws.Range("H4").Value = 78.56999999999999
$ws.Range("L4").ClearContents()

# Row 5 (Pathological_Ear)
$ws.Range("F5").Value = 100
$ws.Range("H5").Value = 100
$ws.Range("I5").Value = 0.2027329192546581
$ws.Range("J5").Value = 0.7536
$ws.Range("K5").ClearContents()

# Row 6 (Left ear)
$ws.Range("F6").Value = 50
$ws.Range("H6").Value = 57.14
$ws.Range("L6").ClearContents()

# Row 7 (Right ear)
$ws.Range("F7").Value = 50
$ws.Range("H7").Value = 42.86
$ws.Range("L7").ClearContents()
